$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The old "Uncovered resources" sub-table used two separate donation
# columns (F "Public/state resources", G "Other funding sources") and
# a third column H computing the remainder (E-F-G). The new layout
# supports a single "Donated" column (F) and computes the remainder in
# column G (E-F), so the old column H is no longer needed.
# ------------------------------------------------------------------
$ws.Columns.Item(8).Delete()

# ---- Table 1 (rows 3-6) header row ----
$ws.Range("F3").Value = "Donated"
$ws.Range("G3").Value = "Total amount of uncovered resources"

# ---- Table 1 data rows ----
$ws.Range("G4").Formula = "=E4-F4"
$ws.Range("G5:G6").Formula = "=E5-F5"

# ---- Table 2 (rows 9-12) header row ----
$ws.Range("F9").Value = "Donated"
$ws.Range("G9").Value = "Total amount of uncovered resources"

# ---- Table 2 data rows ----
$ws.Range("F10").Value = 4500
$ws.Range("G10").Formula = "=E10-F10"
$ws.Range("G11:G12").Formula = "=E11-F11"

# ---- Table 3 (rows 15-18) header row ----
$ws.Range("F15").Value = "Donated"
$ws.Range("G15").Value = "Total amount of uncovered resources"

# ---- Table 3 data rows ----
$ws.Range("G16").Formula = "=E16-F16"
$ws.Range("G17:G18").Formula = "=E17-F17"

# ---- Grand total row ----
$ws.Range("G20").NumberFormat = "#,##0.00"

# Move the active selection like the source workbook (just past the
# bottom of the table, in the new last column).
$ws.Range("G21").Select()
